$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, pushing existing rows 40..123 down to 41..124
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record
$ws.Cells.Item(40, 1).Value = 3
$ws.Cells.Item(40, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(40, 3).Value = "Coquimbo"
$ws.Cells.Item(40, 4).Value = 44581
$ws.Cells.Item(40, 5).Value = 5
$ws.Cells.Item(40, 6).Value = 100112052
$ws.Cells.Item(40, 7).Value = "Albahaca"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 110
$ws.Cells.Item(40, 11).Value = 4000
$ws.Cells.Item(40, 12).Value = 4500
$ws.Cells.Item(40, 13).Value = 4227
$ws.Cells.Item(40, 14).Value = '$/docena de matas'
$ws.Cells.Item(40, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(40, 16).Value = 704
$ws.Cells.Item(40, 17).Value = 6
$ws.Cells.Item(40, 18).Value = "Hortaliza"

# Apply the same date number format used by the rest of column D to the new D40 cell
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
